{"js": "const body = context.document.body;\nbody.clear();\nconst p = body.paragraphs.getFirst();\np.insertText(\"# How to Learn to Draw Manga\\u000bManga is a catch-all term for Japanese comics\u2014a simple, distinctive, and pleasing style that includes a near-infinite array of genres and styles. If you want to learn to draw manga, here are some steps to follow and resources to check out.\\u000b\\u000b## Step 1: Familiarise Yourself with the Style\\u000bManga is an art style as well as a print style, which includes very recognisable visual and storytelling conventions. \\u000b\\u000bStart by understanding these conventions:\\u000b- Manga eyes are bigger than in real life, mouths are smaller, and facial features are simplified.\\u000b- Manga hair often defies gravity.\\u000b- Characters' anatomical proportions are part of what makes manga instantly recognisable.\\u000b- Motion lines can indicate movement or drama.\\u000b- Manga uses visual shortcuts to convey emotions or ideas, such as a giant bead of sweat for frustration or a snot bubble for sleep.\\u000b\\u000bStudy popular manga comics like *Naruto*, *Case Closed*, or *Oh My Goddess!* to understand these visual shortcuts and broad visual styles.\\u000b\\u000b## Step 2: Learn the Fundamentals\\u000bLike any style of drawing, learning the fundamentals of the craft is essential. These include:\\u000b- Anatomy\\u000b- Perspective\\u000b- Composition\\u000b- Light and shadow\\u000b- Proportions\\u000b- Character design\\u000b- Storytelling\\u000b- Inking and line work\\u000b\\u000bOptionally, you could also learn about colour theory, but manga artists draw using mostly black, white, and shades of grey. \\u000b\\u000b## Step 3: Practice Drawing\\u000bAspiring manga artists should practice drawing:\\u000b- Basic shapes like circles and squares, and straight lines.\\u000b- Real-life anatomy, which is a good foundation for adjusting proportions when stylising. Life drawing classes are widely accessible.\\u000b- Manga characters that inspire you, treating the process as a learning experience about the skills used by professional artists. Avoid the temptation to copy or plagiarise their work, but analyse their techniques and practice using your own style.\\u000b\\u000bRemember that drawing is a muscle memory, so the more you practice, the better you'll get. Draw every day, and keep a sketchbook for your exercises. \\u000b\\u000b## Step 4: Use Reference Material\\u000bUse reference material, such as photographs or real objects, to help your practice, especially if you're a beginner. This can help with proportions and placement.\\u000b\\u000b## Step 5: Watch Tutorials and Read Books\\u000bThere are plenty of resources available to help you learn to draw manga, and many are free.\\u000b\\u000b### Websites\\u000b- **AnimeOutline** provides a comprehensive beginner's guide to drawing anime and manga, covering everything from anatomy to shading. There are also lists of common mistakes and how to avoid them, and suggestions for further learning.\\u000b- **How to Draw Manga** is a comprehensive site with all the basics, and plenty of examples.\\u000b- **YouTube** has many channels dedicated to teaching manga and anime art, including **Whyt Manga**, **Manga Materials**, **Sycra**, **Alphonso Dunn**, **Koizu**, **REIQ**, and **Kuzomari**.\\u000b\\u000b### Books\\u000b- **How to Draw Manga: Basics and Beyond!** by Manga University is a good introduction to the art, with more than 1,000 illustrations by Japanese artists.\\u000b- **Manga for the Beginner: Everything you Need to Start Drawing Right Away!** by Christopher Hart details how to draw many different elements, and also covers lighting and shading.\\u000b- **Mastering Manga with Mark Crilley** contains lessons from a master manga artist, including basic drawing lessons and complex art concepts.\\u000b- **Pop Manga: How to Draw the Coolest, Cutest Characters, Animals, Mascots, and More** by Camilla d'Errico and Stephen W. Martin goes beyond a how-to guide with lessons on creating striking characters and capturing action.\\u000b\\u000bThere are many other books available; a simple Amazon search will help you find them.\\u000b\\u000b## Step 6: Develop Your Own Style\\u000b Manga artist Mark Crilley says that the first step to developing your own style is to allow yourself \\\"a period of complete lack of originality\\\". You might start by copying other artists' work to learn their techniques, but eventually, you'll combine these with your own to develop your own style. \\u000b\\u000bSome final tips:\\u000b- Learning to draw is a long process, so be patient and don't give up. \\u000b- Draw your favourite characters to help you decide what style you prefer, but remember not to copy them \u2013 use the process as an opportunity to learn.\\u000b- Don't let others' opinions discourage you. Everyone works at different paces, and you should focus on your own progress.\\u000bGood luck!\", \"Replace\");\nawait context.sync();\n", "ps1": "$texts = @(\n    '# How to Learn to Draw Manga',\n    'Manga is a catch-all term for Japanese comics\u2014a simple, distinctive, and pleasing style that includes a near-infinite array of genres and styles. If you want to learn to draw manga, here are some steps to follow and resources to check out.',\n    '## Step 1: Familiarise Yourself with the Style',\n    'Manga is an art style as well as a print style, which includes very recognisable visual and storytelling conventions. ',\n    'Start by understanding these conventions:',\n    '- Manga eyes are bigger than in real life, mouths are smaller, and facial features are simplified.',\n    '- Manga hair often defies gravity.',\n    '- Characters'' anatomical proportions are part of what makes manga instantly recognisable.',\n    '- Motion lines can indicate movement or drama.',\n    '- Manga uses visual shortcuts to convey emotions or ideas, such as a giant bead of sweat for frustration or a snot bubble for sleep.',\n    'Study popular manga comics like *Naruto*, *Case Closed*, or *Oh My Goddess!* to understand these visual shortcuts and broad visual styles.',\n    '## Step 2: Learn the Fundamentals',\n    'Like any style of drawing, learning the fundamentals of the craft is essential. These include:',\n    '- Anatomy',\n    '- Perspective',\n    '- Composition',\n    '- Light and shadow',\n    '- Proportions',\n    '- Character design',\n    '- Storytelling',\n    '- Inking and line work',\n    'Optionally, you could also learn about colour theory, but manga artists draw using mostly black, white, and shades of grey. ',\n    '## Step 3: Practice Drawing',\n    'Aspiring manga artists should practice drawing:',\n    '- Basic shapes like circles and squares, and straight lines.',\n    '- Real-life anatomy, which is a good foundation for adjusting proportions when stylising. Life drawing classes are widely accessible.',\n    '- Manga characters that inspire you, treating the process as a learning experience about the skills used by professional artists. Avoid the temptation to copy or plagiarise their work, but analyse their techniques and practice using your own style.',\n    'Remember that drawing is a muscle memory, so the more you practice, the better you''ll get. Draw every day, and keep a sketchbook for your exercises. ',\n    '## Step 4: Use Reference Material',\n    'Use reference material, such as photographs or real objects, to help your practice, especially if you''re a beginner. This can help with proportions and placement.',\n    '## Step 5: Watch Tutorials and Read Books',\n    'There are plenty of resources available to help you learn to draw manga, and many are free.',\n    '### Websites',\n    '- **AnimeOutline** provides a comprehensive beginner''s guide to drawing anime and manga, covering everything from anatomy to shading. There are also lists of common mistakes and how to avoid them, and suggestions for further learning.',\n    '- **How to Draw Manga** is a comprehensive site with all the basics, and plenty of examples.',\n    '- **YouTube** has many channels dedicated to teaching manga and anime art, including **Whyt Manga**, **Manga Materials**, **Sycra**, **Alphonso Dunn**, **Koizu**, **REIQ**, and **Kuzomari**.',\n    '### Books',\n    '- **How to Draw Manga: Basics and Beyond!** by Manga University is a good introduction to the art, with more than 1,000 illustrations by Japanese artists.',\n    '- **Manga for the Beginner: Everything you Need to Start Drawing Right Away!** by Christopher Hart details how to draw many different elements, and also covers lighting and shading.',\n    '- **Mastering Manga with Mark Crilley** contains lessons from a master manga artist, including basic drawing lessons and complex art concepts.',\n    '- **Pop Manga: How to Draw the Coolest, Cutest Characters, Animals, Mascots, and More** by Camilla d''Errico and Stephen W. Martin goes beyond a how-to guide with lessons on creating striking characters and capturing action.',\n    'There are many other books available; a simple Amazon search will help you find them.',\n    '## Step 6: Develop Your Own Style',\n    ' Manga artist Mark Crilley says that the first step to developing your own style is to allow yourself \"a period of complete lack of originality\". You might start by copying other artists'' work to learn their techniques, but eventually, you''ll combine these with your own to develop your own style. ',\n    'Some final tips:',\n    '- Learning to draw is a long process, so be patient and don''t give up. ',\n    '- Draw your favourite characters to help you decide what style you prefer, but remember not to copy them \u2013 use the process as an opportunity to learn.',\n    '- Don''t let others'' opinions discourage you. Everyone works at different paces, and you should focus on your own progress.',\n    'Good luck!'\n)\n$breaks = @(1,2,1,2,1,1,1,1,1,2,2,1,1,1,1,1,1,1,1,1,2,2,1,1,1,1,2,2,1,2,1,2,1,1,1,2,1,1,1,1,2,2,1,2,1,1,1,1,0)\n\n$sb = New-Object System.Text.StringBuilder\nfor ($i = 0; $i -lt $texts.Length; $i++) {\n    [void]$sb.Append($texts[$i])\n    for ($j = 0; $j -lt $breaks[$i]; $j++) {\n        [void]$sb.Append([char]11)\n    }\n}\n\n$d = $word.ActiveDocument\n$d.Content.Text = $sb.ToString()\n"}
